$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.360.95'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.923.08'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.8134'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.38'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.97%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9994'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3269'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.24'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07296'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +5.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7940'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +6.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08098'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.927.14'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.418'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +4.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.63'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.388.01'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.39'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.122'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +4.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '252.16'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007893'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.179.92'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.08%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.049'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +17.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1695'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +21.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.559'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.55'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.15'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.170'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +6.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.375'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.549'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.353'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.152'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05657'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.302'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7480'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.87%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.001'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.724'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01960'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.824'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4513'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '74.47'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.997'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8569'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.935'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.038.20'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +5.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9985'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '103.62'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.08%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.72%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.667'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.085.85'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.44%  '
